# Hotfix: Tue Nov  5 14:57:49 RTZ 2024
#
# 1) "Python" sheet: remove the stray test row (id 2110 / "dfgfdgfd").
# 2) "Links" sheet: re-sort rows by id ascending, normalize every URL in
#    column C to use an explicit "https://" scheme, and append a new
#    link (id 51) about working with files/pathlib in Python.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Python sheet — delete row 35 (A35=2110, "dfgfdgfd" / "gdfgfdgfdgdfggfd")
# ---------------------------------------------------------------------
$wsPython = $wb.Worksheets.Item("Python")
$wsPython.Rows.Item(35).Delete()

# ---------------------------------------------------------------------
# 2) Links sheet — sort by column A ascending, fix up URLs, append row
# ---------------------------------------------------------------------
$wsLinks = $wb.Worksheets.Item("Links")

$sortRange = $wsLinks.Range("A1:C15")
$sortRange.Sort($wsLinks.Range("A1:A15"))

for ($r = 1; $r -le 15; $r++) {
    $url = $wsLinks.Cells.Item($r, 3).Value2
    if ($url -notlike "https://*") {
        $wsLinks.Cells.Item($r, 3).Value = "https://" + $url
    }
}

$wsLinks.Cells.Item(16, 1).Value = 51
$wsLinks.Cells.Item(16, 2).Value = "Полезная страница с методами работы с файлами в Python."
$wsLinks.Cells.Item(16, 3).Value = "https://victor-komlev.ru/rabota-s-operatsionnoj-i-fajlovoj-sistemoj-v-python-pathlib-os-shutil/"
